# Updated cryptos list on Fri Dec  8 18:10:31 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value so it is stored as TEXT (matches source data which
# keeps numeric-looking price strings, e.g. '0.671', as literal text) instead
# of letting Excel auto-convert it to a Number/float.
function Set-TextValue($cell, $text) {
    $cell.Formula = '="' + $text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = $false
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.849.83"
$ws.Range("E2").Value = "  +0.91%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.352.33"
$ws.Range("E3").Value = "  +0.61%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.18%  "

# Row 5 - XRP
Set-TextValue $ws.Range("D5") "0.671"
$ws.Range("E5").Value = "  +3.43%  "

# Row 6 - BNB
Set-TextValue $ws.Range("D6") "235.96"
$ws.Range("E6").Value = "  +1.35%  "

# Row 7 - Solana
Set-TextValue $ws.Range("D7") "73.42"
$ws.Range("E7").Value = "  +11.00%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.06%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.545"
$ws.Range("E9").Value = "  +19.74%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.09%  "

# Row 11 - Avalanche
Set-TextValue $ws.Range("D11") "28.18"
$ws.Range("E11").Value = "  +4.69%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +1.94%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.698.89"
$ws.Range("E13").Value = "  +0.58%  "

# Row 14 - Chainlink
Set-TextValue $ws.Range("D14") "16.70"
$ws.Range("E14").Value = "  +7.42%  "

# Row 15 - Polkadot
Set-TextValue $ws.Range("D15") "6.67"
$ws.Range("E15").Value = "  +6.07%  "

# Row 16 - Polygon
Set-TextValue $ws.Range("D16") "0.892"
$ws.Range("E16").Value = "  +4.66%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.352.96"
$ws.Range("E17").Value = "  +1.06%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.808.50"
$ws.Range("E18").Value = "  +0.99%  "

# Row 19 - ShibaInu
Set-TextValue $ws.Range("D19") "0.0000101"
$ws.Range("E19").Value = "  +3.14%  "

# Row 20 - Litecoin
Set-TextValue $ws.Range("D20") "77.03"
$ws.Range("E20").Value = "  +3.61%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  +2.83%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "253.86"
$ws.Range("E22").Value = "  +1.68%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.01%  "

# Row 24 - WEMIXToken
$ws.Range("E24").Value = "  -2.50%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "2.50"
$ws.Range("E25").Value = "  +2.17%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  +6.13%  "

# Row 27 - Toncoin
Set-TextValue $ws.Range("D27") "2.30"
$ws.Range("E27").Value = "  +0.88%  "

# Row 28 - EthereumClassic
Set-TextValue $ws.Range("D28") "22.42"
$ws.Range("E28").Value = "  +1.03%  "

# Row 29 - Monero
Set-TextValue $ws.Range("D29") "172.63"
$ws.Range("E29").Value = "  -1.09%  "

# Row 30 - ImmutableX
Set-TextValue $ws.Range("D30") "1.58"
$ws.Range("E30").Value = "  +6.97%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  +1.78%  "

# Row 32 - Stellar
$ws.Range("E32").Value = "  +5.00%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "5.16"
$ws.Range("E33").Value = "  +2.78%  "

# Row 34 - Hedera
Set-TextValue $ws.Range("D34") "0.0712"
$ws.Range("E34").Value = "  +3.09%  "

# Row 35 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D35") "5.17"
$ws.Range("E35").Value = "  +3.96%  "

# Row 36 - RenderToken
Set-TextValue $ws.Range("D36") "3.92"
$ws.Range("E36").Value = "  +8.22%  "

# Row 37 - LidoDAOToken
$ws.Range("E37").Value = "  -3.60%  "

# Row 38 - THORChain
Set-TextValue $ws.Range("D38") "6.42"
$ws.Range("E38").Value = "  -1.47%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +5.95%  "

# Row 40 - InjectiveProtocol
Set-TextValue $ws.Range("D40") "19.57"
$ws.Range("E40").Value = "  +5.20%  "

# Row 41 - BinanceUSD
$ws.Range("E41").Value = "  -0.09%  "

# Row 42 - FraxShare
Set-TextValue $ws.Range("D42") "8.84"
$ws.Range("E42").Value = "  -2.50%  "

# Row 43 - TrustWalletToken
$ws.Range("E43").Value = "  +1.93%  "

# Row 44 - Cronos
Set-TextValue $ws.Range("D44") "0.0978"
$ws.Range("E44").Value = "  +3.17%  "

# Row 45 - ARBITRUM
$ws.Range("E45").Value = "  -0.94%  "

# Row 46 - FTXToken
$ws.Range("E46").Value = "  +1.64%  "

# Row 47 - Algorand
Set-TextValue $ws.Range("D47") "0.182"
$ws.Range("E47").Value = "  +11.45%  "

# Row 48 - Aave
Set-TextValue $ws.Range("D48") "97.18"
$ws.Range("E48").Value = "  -2.34%  "

# Row 49 - Maker
$ws.Range("D49").Value = "1.436.00"
$ws.Range("E49").Value = "  -0.80%  "

# Row 50 - was HuobiToken, now NEARProtocol
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D50") "2.29"
$ws.Range("E50").Value = "  +0.22%  "

# Row 51 - was RocketPoolETH, now HuobiToken
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D51") "2.78"
$ws.Range("E51").Value = "  +1.29%  "

